$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.514.77"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "2.371.48"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'508.05"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "'134.00"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  -1.62%  "

$ws.Range("D9").Value = "2.393.24"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "'0.0972"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").Value = "'4.84"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("E13").Value = "  -3.57%  "

$ws.Range("D14").Value = "2.799.52"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").Value = "'21.99"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").Value = "56.473.71"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").Value = "2.370.22"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "'10.05"
$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "'312.06"
$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").Value = "'6.29"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "'65.64"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("E25").Value = "  -0.50%  "

$ws.Range("D26").Value = "'0.375"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("E28").Value = "  -0.51%  "

$ws.Range("D29").Value = "'171.17"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").Value = "'5.88"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").Value = "'1.11"

$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  +8.24%  "

$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").Value = "'36.64"
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("D41").Value = "'1.43"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("D42").Value = "'0.377"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.40"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'5.02"
$ws.Range("E44").Value = "  +3.45%  "

$ws.Range("D45").Value = "'127.66"
$ws.Range("E45").Value = "  -3.39%  "

$ws.Range("D46").Value = "'0.564"
$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").Value = "'0.0902"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").Value = "'247.24"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "'0.0487"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").Value = "'17.21"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("D51").Value = "'0.0210"
$ws.Range("E51").Value = "  +0.35%  "
